# Auto-generated Excel COM-interop edit script
# Applies the 'Update automàtic: dades i banners [2026-02-22 18:20]' changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = "2026-02-22 18:18:20"
$ws.Range("E3").Value = "2026-02-22 18:18:22"
$ws.Range("E4").Value = "2026-02-22 18:18:25"
$ws.Range("J4").Value = "1027.2 hPa"
$ws.Range("E5").Value = "2026-02-22 18:18:27"
$ws.Range("E6").Value = "2026-02-22 18:18:30"
$ws.Range("J6").Value = "1027.2 hPa"
$ws.Range("O6").Value = "13.1 °C"
$ws.Range("E7").Value = "2026-02-22 18:18:32"
$ws.Range("E8").Value = "2026-02-22 18:18:34"
$ws.Range("E9").Value = "2026-02-22 18:18:37"
$ws.Range("O9").Value = "11.2 °C"
$ws.Range("E10").Value = "2026-02-22 18:18:40"
$ws.Range("E11").Value = "2026-02-22 18:18:42"
$ws.Range("O11").Value = "8.8 °C"
$ws.Range("E12").Value = "2026-02-22 18:18:44"
$ws.Range("H12").Value = "85%"
$ws.Range("O12").Value = "9.9 °C"
$ws.Range("E13").Value = "2026-02-22 18:18:47"
$ws.Range("J13").Value = "1030.6 hPa"
$ws.Range("O13").Value = "6.2 °C"
$ws.Range("E14").Value = "2026-02-22 18:18:49"
$ws.Range("H14").Value = "70%"
$ws.Range("E15").Value = "2026-02-22 18:18:52"
$ws.Range("E16").Value = "2026-02-22 18:18:54"
$ws.Range("H16").Value = "17%"
$ws.Range("E17").Value = "2026-02-22 18:18:57"
$ws.Range("O17").Value = "10.3 °C"
$ws.Range("E18").Value = "2026-02-22 18:18:59"
$ws.Range("O18").Value = "10.1 °C"
$ws.Range("E19").Value = "2026-02-22 18:19:01"
$ws.Range("H19").Value = "44%"
$ws.Range("E20").Value = "2026-02-22 18:19:04"
$ws.Range("O20").Value = "3.9 °C"
$ws.Range("E21").Value = "2026-02-22 18:19:06"
$ws.Range("H21").Value = "58%"
$ws.Range("J21").Value = "1029.3 hPa"
$ws.Range("O21").Value = "9.1 °C"
$ws.Range("E22").Value = "2026-02-22 18:19:09"
$ws.Range("E23").Value = "2026-02-22 18:19:11"
$ws.Range("L23").Value = "17.3 km/h - 340º 17:35 TU"
$ws.Range("E24").Value = "2026-02-22 18:19:14"
$ws.Range("J24").Value = "1029.8 hPa"
$ws.Range("O24").Value = "7.7 °C"
$ws.Range("E25").Value = "2026-02-22 18:19:16"
$ws.Range("E26").Value = "2026-02-22 18:19:19"
$ws.Range("H26").Value = "33%"
$ws.Range("O26").Value = "11.7 °C"
$ws.Range("E27").Value = "2026-02-22 18:19:21"
$ws.Range("H27").Value = "25%"
$ws.Range("E28").Value = "2026-02-22 18:19:24"
$ws.Range("J28").Value = "1027.7 hPa"
$ws.Range("O28").Value = "10.6 °C"
$ws.Range("E29").Value = "2026-02-22 18:19:26"
$ws.Range("H29").Value = "81%"
$ws.Range("E30").Value = "2026-02-22 18:19:29"
$ws.Range("O30").Value = "12.5 °C"
$ws.Range("E31").Value = "2026-02-22 18:19:31"
$ws.Range("O31").Value = "14.3 °C"
$ws.Range("E32").Value = "2026-02-22 18:19:34"
$ws.Range("O32").Value = "6.5 °C"
$ws.Range("E33").Value = "2026-02-22 18:19:36"
$ws.Range("J33").Value = "1028.8 hPa"
$ws.Range("L33").Value = "13.7 km/h - 112º 17:43 TU"
$ws.Range("O33").Value = "8.2 °C"
$ws.Range("E34").Value = "2026-02-22 18:19:39"
$ws.Range("H34").Value = "43%"
$ws.Range("E35").Value = "2026-02-22 18:19:41"
$ws.Range("J35").Value = "1028.5 hPa"
$ws.Range("E36").Value = "2026-02-22 18:19:44"
$ws.Range("E37").Value = "2026-02-22 18:19:46"
$ws.Range("E38").Value = "2026-02-22 18:19:49"
$ws.Range("H38").Value = "66%"
$ws.Range("O38").Value = "11.6 °C"
$ws.Range("E39").Value = "2026-02-22 18:19:51"
$ws.Range("H39").Value = "26%"
$ws.Range("O39").Value = "5.1 °C"
$ws.Range("E40").Value = "2026-02-22 18:19:54"
$ws.Range("J40").Value = "1029.2 hPa"
$ws.Range("E41").Value = "2026-02-22 18:19:56"
$ws.Range("H41").Value = "74%"
$ws.Range("O41").Value = "11.4 °C"
$ws.Range("E42").Value = "2026-02-22 18:19:58"
$ws.Range("O42").Value = "10.8 °C"
$ws.Range("E43").Value = "2026-02-22 18:20:01"
$ws.Range("H43").Value = "71%"
$ws.Range("O43").Value = "9.1 °C"
$ws.Range("E44").Value = "2026-02-22 18:20:03"
$ws.Range("H44").Value = "39%"
$ws.Range("E45").Value = "2026-02-22 18:20:06"
$ws.Range("K45").Value = "13.8 MJ/m2"
$ws.Range("E46").Value = "2026-02-22 18:20:08"
$ws.Range("J46").Value = "1029.8 hPa"
$ws.Range("O46").Value = "8.9 °C"
